$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.465.52"
$ws.Range("E2").Value = "  -2.61%  "
$ws.Range("D3").Value = "2.484.69"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'313.07"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'94.70"
$ws.Range("E6").Value = "  -4.47%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  -2.66%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -3.71%  "
$ws.Range("D10").Value = "'33.60"
$ws.Range("E10").Value = "  -5.00%  "
$ws.Range("D11").Value = "'0.0781"
$ws.Range("E11").Value = "  -2.61%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").Value = "2.868.60"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "'15.55"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("D16").Value = "2.491.13"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("D18").Value = "41.419.32"
$ws.Range("E18").Value = "  -2.76%  "
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("D20").Value = "0.0₃0923"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "'11.26"
$ws.Range("E21").Value = "  -7.82%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'237.40"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("E24").Value = "  -3.31%  "
$ws.Range("E25").Value = "  -4.70%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'24.18"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").Value = "'9.75"
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("D30").Value = "'36.43"
$ws.Range("D31").Value = "'152.04"
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("E32").Value = "  -5.81%  "
$ws.Range("D33").Value = "'2.57"
$ws.Range("E33").Value = "  -3.80%  "
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'18.11"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.0752"
$ws.Range("E36").Value = "  -4.46%  "
$ws.Range("D37").Value = "'3.09"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").Value = "'0.114"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("E40").Value = "  -7.16%  "
$ws.Range("D41").Value = "'4.19"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").Value = "'1.01"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").Value = "2.007.64"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "'19.55"
$ws.Range("E44").Value = "  -10.27%  "
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("E46").Value = "  -7.22%  "
$ws.Range("D47").Value = "'8.81"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").Value = "2.733.10"
$ws.Range("E48").Value = "  -1.56%  "
$ws.Range("D49").Value = "'69.79"
$ws.Range("D50").Value = "'97.21"
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("D51").Value = "'74.62"
$ws.Range("E51").Value = "  -5.39%  "
